$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B1:B6").Copy($ws.Range("C1:C6"))
$ws.Range("C1").Value = "French"
$ws.Range("B1").ClearContents()
$ws.Range("B2:B6").Clear()
